# Generate Report for Archive
#
# The localization status for every tracked file has moved from
# "Ready for handoff" to "In Translation". Update the Status values on the
# Overview / zh-cn / de-de sheets, then shrink the Status columns so they
# fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: per-language status lives in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = 12.5
